# Add "Manufacturer" and "Manufacturer Product Number" columns to the BOM,
# between the existing "Value" and "Description" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F (existing Description/Datasheet shift to G:H).
$ws.Range("E1:F1").EntireColumn.Insert()

# Header row.
$ws.Range("E1").Value = "Manufacturer"
$ws.Range("F1").Value = "Manufacturer Product Number"

# Row 2 - C1,C2,C3 0.1uF capacitors.
$ws.Range("E2").Value = "KYOCERA AVX"
$ws.Range("F2").Value = "KGM21NR71H104KT"

# Row 3 - J1 14 position header.
$ws.Range("E3").Value = "TE Connectivity AMP Connectors"
$ws.Range("F3").Value = "826632-7"

# Row 4 - J2 9 position header.
$ws.Range("E4").Value = "Adam Tech"
$ws.Range("F4").Value = "PH1-09-UA"

# Row 5 - J3 6 position header.
$ws.Range("E5").Value = "Adam Tech"
$ws.Range("F5").Value = "BHR-06-VUA"

# Row 6 - R1-R12 75k resistors.
$ws.Range("E6").Value = "TE Connectivity Passive Product"
$ws.Range("F6").Value = "CPF0805B75KE1"

# Row 7 - U1,U3 74HC257 multiplexers.
$ws.Range("E7").Value = "Texas Instruments"
$ws.Range("F7").Value = "CD74HC257M96"

# Row 8 - U2 ATTINY404-SSNR microcontroller.
$ws.Range("E8").Value = "Microchip Technology"
$ws.Range("F8").Value = "ATTINY404-SSNR"

# Restore the column widths to the values used by the final layout.
$ws.Columns.Item(5).ColumnWidth = 28.55
$ws.Columns.Item(6).ColumnWidth = 27.58
$ws.Columns.Item(7).ColumnWidth = 87.11
$ws.Columns.Item(8).ColumnWidth = 55.68
